# Add support for JS-evaluated "showif" expressions.
# - Fill the (previously empty) "showif" column E with example formulas
#   for every question row that didn't already have one.
# - Row 15 ("slider_list_longer") gets a distinct example formula that
#   references another field ("abode"); every other touched row gets the
#   generic "mc_religion == 1" example.
# - A handful of rows grow taller to fit the new text, and the sheet
#   selection moves to the newly edited E15 cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

$genericShowIf = "mc_religion == 1"
$specialShowIf = "abode %contains% 'berlin'"

# All the other previously-blank "showif" cells get the generic example.
# (Written before the E15 special-case so the shared-string table gets the
# generic string first, matching the original author's edit order.)
$genericRows = @(8, 9, 10, 11, 12, 14, 16, 17, 18, 19, 20, 21, 22, 23, 24, 25, 26, 27, 28, 37, 38, 40, 41, 42, 43, 44, 45, 46, 47, 48, 49, 50, 51, 52, 53, 54, 55)
foreach ($r in $genericRows) {
    $ws.Range("E$r").Value = $genericShowIf
}

# Row 15 is unique: different example expression.
$ws.Range("E15").Value = $specialShowIf

# Rows that grow to accommodate the newly-added text.
$ws.Rows.Item(15).RowHeight = 45
$ws.Rows.Item(19).RowHeight = 30
$ws.Rows.Item(24).RowHeight = 30
$ws.Rows.Item(25).RowHeight = 30
$ws.Rows.Item(26).RowHeight = 34
$ws.Rows.Item(40).RowHeight = 30
$ws.Rows.Item(41).RowHeight = 30
$ws.Rows.Item(44).RowHeight = 30
$ws.Rows.Item(46).RowHeight = 30
$ws.Rows.Item(49).RowHeight = 30
$ws.Rows.Item(50).RowHeight = 30
$ws.Rows.Item(51).RowHeight = 30

# Move the visible selection to the cell that was just edited.
$ws.Activate()
$ws.Range("E15").Select()
